$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Add new cell D3 with value "LMS" (new shared string)
$ws.Range("D3").Value = "LMS"

# Update the selection on the Login sheet to C3
$ws.Range("C3").Select()
